$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dRange = $ws.Range("D2:D48")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "89.773.15"
$ws.Range("E2").Value = "  +3.03%  "
$ws.Range("D3").Value = "3.299.64"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "213.44"
$ws.Range("E5").Value = "  -3.13%  "
$ws.Range("D6").Value = "630.27"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("D7").Value = "0.386"
$ws.Range("E7").Value = "  +17.59%  "
$ws.Range("D8").Value = "0.712"
$ws.Range("E8").Value = "  +13.91%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("E11").Value = "  -5.15%  "
$ws.Range("D12").Value = "0.190"
$ws.Range("E12").Value = "  +13.41%  "
$ws.Range("E13").Value = "  -4.26%  "
$ws.Range("D14").Value = "34.42"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "3.907.96"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").Value = "5.47"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "89.428.14"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").Value = "3.310.85"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("E19").Value = "  -4.06%  "
$ws.Range("D20").Value = "3.10"
$ws.Range("E20").Value = "  -3.99%  "
$ws.Range("D21").Value = "439.70"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("D22").Value = "8.95"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").Value = "7.47"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "5.26"
$ws.Range("E25").Value = "  -2.96%  "
$ws.Range("D26").Value = "12.24"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").Value = "3.464.60"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  +3.61%  "
$ws.Range("D29").Value = "77.29"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "0.185"
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "8.90"
$ws.Range("D34").Value = "564.20"
$ws.Range("E34").Value = "  -7.11%  "
$ws.Range("E35").Value = "  -11.33%  "
$ws.Range("E36").Value = "  +9.15%  "
$ws.Range("E37").Value = "  -4.57%  "
$ws.Range("D38").Value = "0.143"
$ws.Range("E38").Value = "  -6.27%  "
$ws.Range("D39").Value = "22.87"
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("D40").Value = "21.85"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").Value = "0.405"
$ws.Range("E43").Value = "  -4.65%  "
$ws.Range("D44").Value = "2.05"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "153.74"
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("D47").Value = "182.06"
$ws.Range("E47").Value = "  -4.82%  "
$ws.Range("D48").Value = "45.04"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("E49").Value = "  +16.00%  "
$ws.Range("E50").Value = "  -4.22%  "
$ws.Range("E51").Value = "  -1.64%  "

$dRange.ClearFormats()
